# Merge the "FWPolicy" sheet into "Global" (renamed to "Sheet1") and add
# a default outbound firewall-policy row, per the commit:
#   "WIP - including default outbound firewall policy"

$wb = $excel.ActiveWorkbook

# --- 1. Drop the separate FWPolicy sheet; it's being folded into the
#        single remaining sheet as extra columns instead. -------------
$fw = $wb.Worksheets.Item("FWPolicy")
$fw.Delete()

# --- 2. Rename the remaining "Global" sheet to "Sheet1". --------------
$ws = $wb.Worksheets.Item("Global")
$ws.Name = "Sheet1"

# --- 3. Add the firewall-policy header columns (U:AD). ----------------
$ws.Range("U1").Value = "policyid"
$ws.Range("V1").Value = "name"
$ws.Range("W1").Value = "source_int"
$ws.Range("X1").Value = "dest_int"
$ws.Range("Y1").Value = "source_add"
$ws.Range("Z1").Value = "dest_add"
$ws.Range("AA1").Value = "service"
$ws.Range("AB1").Value = "action"
$ws.Range("AC1").Value = "schedule"
$ws.Range("AD1").Value = "logtraffic"

# --- 4. Populate the default outbound policy on row 2 (HQ). -----------
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = "Default-Outbound"
$ws.Range("V2").WrapText = $true
$ws.Range("W2").Value = "port1"
$ws.Range("X2").Value = "port2"
$ws.Range("Y2").Value = "all"
$ws.Range("Z2").Value = "all"
$ws.Range("AA2").Value = "ALL"
$ws.Range("AB2").Value = "accept"
$ws.Range("AC2").Value = "always"
$ws.Range("AD2").Value = "all"
